$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J41").Value = 367.2
$ws.Range("L41").Value = 367.2
$ws.Range("M41").Value = -606.1875
$ws.Range("N41").Value = -1247.2
$ws.Range("I41").Value = 1046.1875
$ws.Range("K41").Value = 1046.1875
$ws.Range("H41").Value = 884.5238000000001
$ws.Range("J53").Value = 345.83334
$ws.Range("H53").Value = 446.75
$ws.Range("L53").Value = 345.83334
$ws.Range("I53").Value = 749.5
$ws.Range("N53").Value = -1619.83334
$ws.Range("K53").Value = 749.5
$ws.Range("M53").Value = -112.5
$ws.Range("K74").Value = 3320.1428
$ws.Range("I74").Value = 3320.1428
$ws.Range("M74").Value = -2384.1428
$ws.Range("H74").Value = 4403.154
$ws.Range("K77").Value = 16600.714
$ws.Range("I77").Value = 3320.1428
$ws.Range("M77").Value = -11920.714
$ws.Range("H77").Value = 4403.154
$ws.Range("N111").Value = -8052.5
$ws.Range("K111").Value = 1014.375
$ws.Range("M111").Value = 2052.625
$ws.Range("I111").Value = 338.125
$ws.Range("H111").Value = 398.4
$ws.Range("J111").Value = 639.5
$ws.Range("L111").Value = 1918.5
$ws.Range("L134").Value = 107000
$ws.Range("J134").Value = 107000
$ws.Range("N134").Value = -117140
$ws.Range("H134").Value = 107000
$ws.Range("J137").Value = 1342.5
$ws.Range("K137").Value = 3219.5295
$ws.Range("M137").Value = -669.5295000000001
$ws.Range("I137").Value = 1073.1765
$ws.Range("H137").Value = 1143.4348
$ws.Range("N137").Value = -9127.5
$ws.Range("L137").Value = 4027.5
$ws.Range("I141").Value = 2528.8
$ws.Range("N141").Value = -14560
$ws.Range("J141").Value = 1400
$ws.Range("H141").Value = 2426.182
$ws.Range("M141").Value = -2406.400000000001
$ws.Range("K141").Value = 7586.400000000001
$ws.Range("L141").Value = 4200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 4622.5557
$ws.Range("H32").Value = 5229.4375
$ws.Range("L32").Value = 14332.667
$ws.Range("N32").Value = -14906.667
$ws.Range("M32").Value = -4335.5557
$ws.Range("K32").Value = 4622.5557
$ws.Range("J32").Value = 14332.667
$ws.Range("I61").Value = 1704.25
$ws.Range("K61").Value = 1704.25
$ws.Range("H61").Value = 1704.25
$ws.Range("M61").Value = -1492.25
$ws.Range("N63").Value = -3375
$ws.Range("L63").Value = 2003
$ws.Range("J63").Value = 2003
$ws.Range("H63").Value = 2376.5
$ws.Range("L66").Value = 10015
$ws.Range("J66").Value = 2003
$ws.Range("H66").Value = 2376.5
$ws.Range("N66").Value = -16879
$ws.Range("K74").Value = 1954.8572
$ws.Range("I74").Value = 1954.8572
$ws.Range("J74").Value = 3259.5
$ws.Range("M74").Value = -1080.8572
$ws.Range("L74").Value = 3259.5
$ws.Range("H74").Value = 2244.7778
$ws.Range("N74").Value = -5007.5
$ws.Range("J77").Value = 3259.5
$ws.Range("K77").Value = 9774.286
$ws.Range("N77").Value = -25033.5
$ws.Range("L77").Value = 16297.5
$ws.Range("I77").Value = 1954.8572
$ws.Range("M77").Value = -5406.286
$ws.Range("H77").Value = 2244.7778
$ws.Range("H136").Value = 1704.25
$ws.Range("I136").Value = 1704.25
$ws.Range("K136").Value = 5112.75
$ws.Range("M136").Value = -2562.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N82").Value = -25766
$ws.Range("J82").Value = 25000
$ws.Range("H82").Value = 21616.5
$ws.Range("L82").Value = 25000
$ws.Range("J85").Value = 25000
$ws.Range("N85").Value = -27652
$ws.Range("H85").Value = 21616.5
$ws.Range("L85").Value = 25000
$ws.Range("K86").Value = 2414.5356
$ws.Range("N86").Value = -4279
$ws.Range("L86").Value = 2033
$ws.Range("I86").Value = 2414.5356
$ws.Range("J86").Value = 2033
$ws.Range("H86").Value = 2377.6128
$ws.Range("M86").Value = -1291.5356
$ws.Range("J89").Value = 2033
$ws.Range("I89").Value = 2414.5356
$ws.Range("M89").Value = -6456.678
$ws.Range("H89").Value = 2377.6128
$ws.Range("K89").Value = 12072.678
$ws.Range("L89").Value = 10165
$ws.Range("N89").Value = -21397

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -724.7143
$ws.Range("K16").Value = 1011.7143
$ws.Range("H16").Value = 1250.5
$ws.Range("I16").Value = 1011.7143
$ws.Range("K86").Value = 0
$ws.Range("N86").Value = -7245
$ws.Range("L86").Value = 4999
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4999
$ws.Range("H86").Value = 4999
$ws.Range("M86").ClearContents()
$ws.Range("J89").Value = 4999
$ws.Range("I89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H89").Value = 4999
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 24995
$ws.Range("N89").Value = -36227
$ws.Range("L99").Value = 2875.3
$ws.Range("M99").Value = -10673.167
$ws.Range("N99").Value = -5871.3
$ws.Range("J99").Value = 2875.3
$ws.Range("K99").Value = 12171.167
$ws.Range("H99").Value = 7945.773
$ws.Range("I99").Value = 12171.167
$ws.Range("H107").Value = 1697.6296
$ws.Range("I107").Value = 1294.4615
$ws.Range("N107").Value = -5912
$ws.Range("M107").Value = 625.5385000000001
$ws.Range("K107").Value = 1294.4615
$ws.Range("L107").Value = 2072
$ws.Range("J107").Value = 2072
$ws.Range("K113").Value = 1011.7143
$ws.Range("M113").Value = 1158.2857
$ws.Range("H113").Value = 1250.5
$ws.Range("I113").Value = 1011.7143
$ws.Range("H126").Value = 7945.773
$ws.Range("L126").Value = 8625.900000000001
$ws.Range("I126").Value = 12171.167
$ws.Range("N126").Value = -13565.9
$ws.Range("K126").Value = 36513.501
$ws.Range("J126").Value = 2875.3
$ws.Range("M126").Value = -34043.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 1250329.5
$ws.Range("M4").Value = -3750876.5
$ws.Range("K4").Value = 3750988.5
$ws.Range("L4").Value = 10463360.4
$ws.Range("N4").Value = -10463584.4
$ws.Range("H4").Value = 1830411.1
$ws.Range("J4").Value = 3487786.8
$ws.Range("H23").Value = 106.4
$ws.Range("J23").Value = 89.375
$ws.Range("N23").Value = -738.125
$ws.Range("L23").Value = 268.125
$ws.Range("J137").Value = 3440.8
$ws.Range("K137").Value = 20004040.5
$ws.Range("M137").Value = -19998940.5
$ws.Range("I137").Value = 6668013.5
$ws.Range("H137").Value = 2859686.2
$ws.Range("N137").Value = -20522.4
$ws.Range("L137").Value = 10322.4
$ws.Range("K140").Value = 4284.5454
$ws.Range("M140").Value = 895.4546
$ws.Range("I140").Value = 1428.1818
$ws.Range("H140").Value = 1750.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("I97").Value = 30303.04
$ws.Range("J97").Value = 2292.6667
$ws.Range("L97").Value = 2292.6667
$ws.Range("H97").Value = 25051.094
$ws.Range("K97").Value = 30303.04
$ws.Range("M97").Value = -29807.04
$ws.Range("N97").Value = -3284.6667
$ws.Range("J132").Value = 959
$ws.Range("N132").Value = -7937
$ws.Range("H132").Value = 3292.9092
$ws.Range("I132").Value = 4626.5713
$ws.Range("L132").Value = 2877
$ws.Range("M132").Value = -11349.7139
$ws.Range("K132").Value = 13879.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14106.781
$ws.Range("M7").Value = -33429.7
$ws.Range("K7").Value = 33541.7
$ws.Range("J7").Value = 5272.727
$ws.Range("L7").Value = 5272.727
$ws.Range("N7").Value = -5496.727
$ws.Range("I7").Value = 33541.7
$ws.Range("K22").Value = 3850.2188
$ws.Range("M22").Value = -3555.2188
$ws.Range("L22").Value = 5110.25
$ws.Range("N22").Value = -5700.25
$ws.Range("I22").Value = 3850.2188
$ws.Range("H22").Value = 3990.2222
$ws.Range("J22").Value = 5110.25
$ws.Range("L27").Value = 5110.25
$ws.Range("K27").Value = 3850.2188
$ws.Range("J27").Value = 5110.25
$ws.Range("N27").Value = -5324.25
$ws.Range("I27").Value = 3850.2188
$ws.Range("M27").Value = -3743.2188
$ws.Range("H27").Value = 3990.2222
$ws.Range("K40").Value = 18678.143
$ws.Range("M40").Value = -18542.143
$ws.Range("H40").Value = 13096.5
$ws.Range("I40").Value = 18678.143
$ws.Range("M82").Value = -913.0999999999999
$ws.Range("I82").Value = 1274.1
$ws.Range("K82").Value = 1274.1
$ws.Range("H82").Value = 1274.1
$ws.Range("I85").Value = 1274.1
$ws.Range("K85").Value = 1274.1
$ws.Range("H85").Value = 1274.1
$ws.Range("M85").Value = -26.09999999999991
$ws.Range("H122").Value = 4051.4443
$ws.Range("M122").Value = -6615.000100000001
$ws.Range("K122").Value = 9065.000100000001
$ws.Range("I122").Value = 3021.6667
$ws.Range("H126").Value = 14106.781
$ws.Range("L126").Value = 15818.181
$ws.Range("I126").Value = 33541.7
$ws.Range("N126").Value = -20758.181
$ws.Range("K126").Value = 100625.1
$ws.Range("J126").Value = 5272.727
$ws.Range("M126").Value = -98155.09999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 14388.889
$ws.Range("K2").Value = 14388.889
$ws.Range("M2").Value = -14276.889
$ws.Range("I2").Value = 14388.889
$ws.Range("L46").Value = 59480.5
$ws.Range("J46").Value = 59480.5
$ws.Range("N46").Value = -59942.5
$ws.Range("H46").Value = 59480.5
$ws.Range("I81").Value = 1640.2858
$ws.Range("J81").Value = 700
$ws.Range("L81").Value = 1400
$ws.Range("K81").Value = 3280.5716
$ws.Range("H81").Value = 1522.75
$ws.Range("N81").Value = -3522
$ws.Range("M81").Value = -2219.5716
$ws.Range("L84").Value = 7000
$ws.Range("N84").Value = -17608
$ws.Range("H84").Value = 1522.75
$ws.Range("M84").Value = -11098.858
$ws.Range("I84").Value = 1640.2858
$ws.Range("J84").Value = 700
$ws.Range("K84").Value = 16402.858
$ws.Range("H126").Value = 1691.1818
$ws.Range("L126").Value = 5998.5
$ws.Range("I126").Value = 1660.35
$ws.Range("N126").Value = -10938.5
$ws.Range("K126").Value = 4981.049999999999
$ws.Range("J126").Value = 1999.5
$ws.Range("M126").Value = -2511.049999999999
$ws.Range("H132").Value = 2007.7556
$ws.Range("I132").Value = 2007.7556
$ws.Range("M132").Value = -3493.266799999999
$ws.Range("K132").Value = 6023.266799999999
$ws.Range("L134").Value = 178441.5
$ws.Range("J134").Value = 59480.5
$ws.Range("N134").Value = -183511.5
$ws.Range("H134").Value = 59480.5

